# Add the WPI reference-electrode URL as a hyperlink in cell E13.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bioreactors")

$url = "https://www.wpiinc.com/driref-2-dri-ref-reference-electrode-2-mm.html?srsltid=AfmBOoo0TIbQbkCCOFrlG2X3Q0pexQVHKHyk5CYPa4OnMnZRMZF5dwtG"

# Insert the hyperlink; Excel will write the URL text into the cell and
# apply the built-in "Hyperlink" style (underline, theme color).
$ws.Hyperlinks.Add($ws.Range("E13"), $url)

# Widen column E so the long URL is visible, matching the manual resize
# that accompanied the new link.
$ws.Columns.Item(5).ColumnWidth = 39

# Reflect the cell the author had selected after making the edit.
$ws.Range("E17").Select()

$wb.Save()
